$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "73.417.36"
$ws.Range("E2").Value = "  -0.16%  "

# Row 3
$ws.Range("D3").Value = "3.986.32"
$ws.Range("E3").Value = "  -1.84%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.12%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "610.52"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +6.15%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "172.13"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +13.19%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.686"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.57%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.04%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.794"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.54%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.187"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +8.67%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "57.05"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.78%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000338"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.37%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.52"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.48%  "

# Row 14
$ws.Range("D14").Value = "4.625.11"
$ws.Range("E14").Value = "  -1.91%  "

# Row 15
$ws.Range("D15").Value = "3.990.40"
$ws.Range("E15").Value = "  -1.89%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.32"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.78%  "

# Row 17
$ws.Range("B17").Value = "Polygon"
$ws.Range("C17").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.24"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.79%  "

# Row 18
$ws.Range("B18").Value = "Chainlink"
$ws.Range("C18").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "20.90"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.28%  "

# Row 19
$ws.Range("D19").Value = "73.342.44"
$ws.Range("E19").Value = "  -0.06%  "

# Row 20
$ws.Range("E20").Value = "  -1.18%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "462.97"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.61%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.86"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +5.37%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "96.53"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.14%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.41"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -5.82%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "14.34"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.88%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.25"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.07%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.14"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.40%  "

# Row 28
$ws.Range("B28").Value = "Filecoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.59"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.83%  "

# Row 29
$ws.Range("B29").Value = "LEO"
$ws.Range("C29").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.93"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.60%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "36.45"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.39%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.94"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.83%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "13.98"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.35%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0000106"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +16.15%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.130"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.36%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "48.54"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.41%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "70.37"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.33%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "638.12"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -7.47%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.433"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.35%  "

# Row 39
$ws.Range("B39").Value = "ThetaToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.40"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.26%  "

# Row 40
$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.148"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.21%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.14%  "

# Row 42
$ws.Range("E42").Value = "  -0.02%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.28"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +38.83%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0485"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.67%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.65"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -5.57%  "

# Row 46
$ws.Range("B46").Value = "FLOKI"
$ws.Range("C46").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.000308"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +10.87%  "

# Row 47
$ws.Range("B47").Value = "Stellar"
$ws.Range("C47").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.149"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.27%  "

# Row 48
$ws.Range("B48").Value = "ApeXProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.47"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.36%  "

# Row 49
$ws.Range("B49").Value = "Fetch.AI"
$ws.Range("C49").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.60"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -5.38%  "

# Row 50
$ws.Range("B50").Value = "WEMIXToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.82"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -15.64%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.01"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.00%  "
